$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out column B rows 2-4 (no longer used)
$ws.Range("B2:B4").ClearContents()

# Update remaining values per the new layout (column A first, then header)
$ws.Range("A2").Value = "https://google.com"
$ws.Range("A3").Value = "youtube.com"
$ws.Range("A4").Value = "google.com/asdf"
$ws.Range("B1").Value = "Response Code"

# Update the active selection shown in the sheet view
$ws.Range("M12").Select()
